$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.114.85"
$ws.Range("E2").Value = "  +1.63%  "

$ws.Range("D3").Value = "2.356.60"
$ws.Range("E3").Value = "  -0.34%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.679"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.93%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.26%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.76"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +6.30%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.566"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +23.60%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.102"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.56%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "32.15"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +20.96%  "

$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +18.89%  "

$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.108"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.77%  "

$ws.Range("D14").Value = "2.704.98"
$ws.Range("E14").Value = "  -0.76%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.89"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.64%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.920"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.88%  "

$ws.Range("D17").Value = "2.353.68"
$ws.Range("E17").Value = "  -0.44%  "

$ws.Range("D18").Value = "44.216.67"
$ws.Range("E18").Value = "  +1.91%  "

$ws.Range("E19").Value = "  +3.98%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "78.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.16%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "256.22"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.05%  "

$ws.Range("E23").Value = "  -0.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.88%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.38%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.80"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.46%  "

$ws.Range("E27").Value = "  +3.96%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.79"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.46%  "

$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "174.81"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.58%  "

$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.59"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.52%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.132"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.88%  "

$ws.Range("E32").Value = "  +5.26%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.38"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.80%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0748"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.21%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.30"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.55%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.88"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.51%  "

$ws.Range("E37").Value = "  -0.41%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.55"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.97%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0273"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.28%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.22"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.42%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.98"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.23%  "

$ws.Range("E42").Value = "  -0.14%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.26"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.61%  "

$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.100"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.03%  "

$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +10.45%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.67%  "

$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.21%  "

$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.186"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +12.03%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.47"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.48%  "

$ws.Range("D50").Value = "1.452.85"
$ws.Range("E50").Value = "  +0.39%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000205"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.61%  "
